$wb = $excel.ActiveWorkbook

# --- Overview sheet: update "Latest HO Xliff Generate Date" for the
#     ce35765b-08c0-4403-b1ec-e01e4fb0289b.md row (row 4, column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-12 16:56:49"

# --- zh-cn sheet: update "Correspond Handoff Datetime" (H4) and
#     "Correspond Handback DateTime" (K4) for the ce35765b row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-12 16:56:42"
$wsZhCn.Range("K4").Value = "2016-08-12 16:57:19"

# --- de-de sheet: update "Correspond Handback DateTime" (K4) for the
#     ce35765b row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-12 16:57:28"
